$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new daily row (row 50) with the next day's data
$ws.Range("A50").Value = 45635
$ws.Range("B50").Value = 127
$ws.Range("C50").Value = 113
$ws.Range("D50").Value = 119

# The new last row (A50) takes on the "date only" number format that used
# to belong to the previous last row.
$ws.Range("A50").NumberFormat = "YYYY-MM-DD"

# The previous last row (A49) becomes a regular data row, so it switches to
# the "date + time" number format shared by all the other data rows.
$ws.Range("A49").NumberFormat = "YYYY-MM-DD HH:MM:SS"
